# Updated Reviews Feed and Reviews Insights
# - Insert a new "Response_Status" sheet (Status / Pending Approval / Rejected / All)
#   right before the "Global Filters" sheet.
# - The new sheet becomes the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Insert new worksheet directly before "Global Filters" so it lands in the
# same tab position the diff shows (index 8, sheetId 9, rId8).
$globalFilters = $wb.Worksheets.Item("Global Filters")
$ws = $wb.Worksheets.Add($globalFilters)
$ws.Name = "Response_Status"

# Populate the data.
$ws.Range("A1").Value = "Status"
$ws.Range("A2").Value = "Pending Approval"
$ws.Range("A3").Value = "Rejected"
$ws.Range("A4").Value = "All"

# Header styling (bold), matching the other lookup-table sheets in the workbook.
$font = $ws.Range("A1").Font
$font.Bold = $true

# Match default print setup used elsewhere in the workbook.
$pageSetup = $ws.PageSetup
$pageSetup.Orientation = 1
$pageSetup.PaperSize = 9

# Leave the selection on the last populated cell, as in the authored workbook.
$ws.Range("A4").Select() | Out-Null
